$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency "Price" column (D) holds numeric-looking strings (e.g.
# "1.00", "0.0000137", "59.360.83") that must stay plain text exactly as
# scraped. Prepending an apostrophe, like typing it into Excel by hand,
# forces a text entry instead of Excel auto-converting it to a number.
$apos = "'"

$ws.Range('D2').Value = $apos + '59.360.83'
$ws.Range('E2').Value = '  -0.30%  '
$ws.Range('D3').Value = $apos + '2.581.27'
$ws.Range('E3').Value = '  -0.35%  '
$ws.Range('D4').Value = $apos + '1.00'
$ws.Range('E4').Value = '  -6.56%  '
$ws.Range('D5').Value = $apos + '552.41'
$ws.Range('E5').Value = '  -2.53%  '
$ws.Range('D6').Value = $apos + '140.13'
$ws.Range('E6').Value = '  -1.94%  '
$ws.Range('D7').Value = $apos + '0.999'
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = $apos + '0.593'
$ws.Range('E8').Value = '  -0.29%  '
$ws.Range('D9').Value = $apos + '2.594.96'
$ws.Range('E9').Value = '  -0.24%  '
$ws.Range('E10').Value = '  +1.43%  '
$ws.Range('D11').Value = $apos + '0.104'
$ws.Range('E11').Value = '  +0.91%  '
$ws.Range('E12').Value = '  +5.44%  '
$ws.Range('D13').Value = $apos + '0.354'
$ws.Range('E13').Value = '  +3.95%  '
$ws.Range('D14').Value = $apos + '3.039.66'
$ws.Range('E14').Value = '  -0.30%  '
$ws.Range('D15').Value = $apos + '59.379.60'
$ws.Range('E15').Value = '  -0.38%  '
$ws.Range('D16').Value = $apos + '23.08'
$ws.Range('E16').Value = '  +5.74%  '
$ws.Range('D17').Value = $apos + '0.0000137'
$ws.Range('E17').Value = '  +0.57%  '
$ws.Range('D18').Value = $apos + '2.597.87'
$ws.Range('D19').Value = $apos + '4.54'
$ws.Range('E19').Value = '  +0.53%  '
$ws.Range('D20').Value = $apos + '338.91'
$ws.Range('E20').Value = '  +0.64%  '
$ws.Range('D21').Value = $apos + '10.35'
$ws.Range('E21').Value = '  +1.48%  '
$ws.Range('D22').Value = $apos + '6.48'
$ws.Range('E22').Value = '  +3.99%  '
$ws.Range('D23').Value = $apos + '0.999'
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').Value = $apos + '0.478'
$ws.Range('E24').Value = '  +7.59%  '
$ws.Range('D25').Value = $apos + '63.17'
$ws.Range('E25').Value = '  -3.02%  '
$ws.Range('D26').Value = $apos + '0.997'
$ws.Range('E26').Value = '  -0.24%  '
$ws.Range('D27').Value = $apos + '0.159'
$ws.Range('E27').Value = '  -1.49%  '
$ws.Range('D28').Value = $apos + '7.46'
$ws.Range('E28').Value = '  +2.84%  '
$ws.Range('D29').Value = $apos + '0.0₃0771'
$ws.Range('E29').Value = '  -1.14%  '
$ws.Range('D30').Value = $apos + '0.998'
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('D31').Value = $apos + '1.68'
$ws.Range('E31').Value = '  -0.39%  '
$ws.Range('D32').Value = $apos + '6.12'
$ws.Range('E32').Value = '  +1.52%  '
$ws.Range('D33').Value = $apos + '157.45'
$ws.Range('E33').Value = '  -1.78%  '
$ws.Range('D34').Value = $apos + '19.11'
$ws.Range('E34').Value = '  +0.82%  '
$ws.Range('D35').Value = $apos + '4.12'
$ws.Range('E35').Value = '  +2.47%  '
$ws.Range('D36').Value = $apos + '1.16'
$ws.Range('E36').Value = '  +2.79%  '
$ws.Range('D37').Value = $apos + '0.898'
$ws.Range('E37').Value = '  +1.11%  '
$ws.Range('D38').Value = $apos + '37.56'
$ws.Range('E38').Value = '  +1.47%  '
$ws.Range('D39').Value = $apos + '1.47'
$ws.Range('E39').Value = '  -0.74%  '
$ws.Range('D40').Value = $apos + '0.839'
$ws.Range('E40').Value = '  -4.29%  '
$ws.Range('D41').Value = $apos + '3.66'
$ws.Range('E41').Value = '  +0.95%  '
$ws.Range('D42').Value = $apos + '288.53'
$ws.Range('E42').Value = '  -1.11%  '
$ws.Range('D43').Value = $apos + '135.74'
$ws.Range('E43').Value = '  +8.95%  '
$ws.Range('D44').Value = $apos + '0.999'
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').Value = $apos + '0.0970'
$ws.Range('E45').Value = '  -0.50%  '
$ws.Range('D46').Value = $apos + '0.597'
$ws.Range('E46').Value = '  +0.71%  '
$ws.Range('D47').Value = $apos + '10.69'
$ws.Range('E47').Value = '  +0.30%  '
$ws.Range('D48').Value = $apos + '0.0532'
$ws.Range('E48').Value = '  -0.62%  '
$ws.Range('D49').Value = $apos + '0.0234'
$ws.Range('E49').Value = '  +1.02%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = $apos + '1.980.11'
$ws.Range('E50').Value = '  +2.60%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = $apos + '18.63'
$ws.Range('E51').Value = '  +1.54%  '
